$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    if ($ws.Range("F2").Value2 -eq 266) { $ws.Range("F2").Value2 = 268 }
    if ($ws.Range("F5").Value2 -eq 6658) { $ws.Range("F5").Value2 = 6666 }
    if ($ws.Range("F6").Value2 -eq 5435) { $ws.Range("F6").Value2 = 5447 }
    if ($ws.Range("F7").Value2 -eq 450) { $ws.Range("F7").Value2 = 451 }
}

$ws1 = $wb.Worksheets.Item("展览")
if ($ws1.Range("F12").Value2 -eq 79) { $ws1.Range("F12").Value2 = 100 }

$ws4 = $wb.Worksheets.Item("全部类型")
if ($ws4.Range("F14").Value2 -eq 80) { $ws4.Range("F14").Value2 = 100 }
